$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 12.142858
$ws.Range("I6").Value = 12.142858
$ws.Range("K6").Value = 36.428574
$ws.Range("M6").Value = 75.571426
$ws.Range("H32").Value = 869
$ws.Range("J32").Value = 956.8
$ws.Range("L32").Value = 956.8
$ws.Range("N32").Value = -1608.8
$ws.Range("H40").Value = 2515.6
$ws.Range("I40").Value = 2515.6
$ws.Range("K40").Value = 2515.6
$ws.Range("M40").Value = -2340.6
$ws.Range("H74").Value = 8375.5
$ws.Range("I74").Value = 6001.25
$ws.Range("J74").Value = 9958.333000000001
$ws.Range("K74").Value = 6001.25
$ws.Range("L74").Value = 9958.333000000001
$ws.Range("M74").Value = -5065.25
$ws.Range("N74").Value = -11830.333
$ws.Range("H77").Value = 8375.5
$ws.Range("I77").Value = 6001.25
$ws.Range("J77").Value = 9958.333000000001
$ws.Range("K77").Value = 30006.25
$ws.Range("L77").Value = 49791.665
$ws.Range("M77").Value = -25326.25
$ws.Range("N77").Value = -59151.665
$ws.Range("H100").Value = 3130.3333
$ws.Range("I100").Value = 2699.5
$ws.Range("K100").Value = 2699.5
$ws.Range("M100").Value = -2158.5
$ws.Range("H113").Value = 5199.5
$ws.Range("I113").Value = 5199.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 5199.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1945.5
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 1447.56
$ws.Range("I132").Value = 1185.9546
$ws.Range("J132").Value = 3366
$ws.Range("K132").Value = 3557.8638
$ws.Range("L132").Value = 10098
$ws.Range("M132").Value = -1027.8638
$ws.Range("N132").Value = -15158
$ws.Range("H138").Value = 2753.7073
$ws.Range("I138").Value = 1973.2727
$ws.Range("J138").Value = 3657.3684
$ws.Range("K138").Value = 5919.8181
$ws.Range("L138").Value = 10972.1052
$ws.Range("M138").Value = -779.8181000000004
$ws.Range("N138").Value = -21252.1052

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2988.45
$ws.Range("I88").Value = 1808
$ws.Range("K88").Value = 1808
$ws.Range("M88").Value = -1402
$ws.Range("H91").Value = 2988.45
$ws.Range("I91").Value = 1808
$ws.Range("K91").Value = 1808
$ws.Range("M91").Value = -404
$ws.Range("H97").Value = 69330
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 69330
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 69330
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -70322
$ws.Range("H122").Value = 5930.393
$ws.Range("I122").Value = 6027.364
$ws.Range("K122").Value = 18082.092
$ws.Range("M122").Value = -15632.092

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2126.7896
$ws.Range("I86").Value = 1775.5625
$ws.Range("K86").Value = 1775.5625
$ws.Range("M86").Value = -652.5625
$ws.Range("H89").Value = 2126.7896
$ws.Range("I89").Value = 1775.5625
$ws.Range("K89").Value = 8877.8125
$ws.Range("M89").Value = -3261.8125
$ws.Range("H94").Value = 1966.6666
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 1966.6666
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 1966.6666
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -2868.6666

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 69000
$ws.Range("J20").Value = 69000
$ws.Range("L20").Value = 69000
$ws.Range("N20").Value = -69472
$ws.Range("H30").Value = 69000
$ws.Range("J30").Value = 69000
$ws.Range("L30").Value = 69000
$ws.Range("N30").Value = -69182
$ws.Range("H31").Value = 4311.2354
$ws.Range("J31").Value = 6289.0527
$ws.Range("L31").Value = 6289.0527
$ws.Range("N31").Value = -6879.0527
$ws.Range("H34").Value = 4311.2354
$ws.Range("J34").Value = 6289.0527
$ws.Range("L34").Value = 6289.0527
$ws.Range("N34").Value = -6693.0527
$ws.Range("H122").Value = 2973.5557
$ws.Range("I122").Value = 3008.8572
$ws.Range("J122").Value = 2850
$ws.Range("K122").Value = 9026.571599999999
$ws.Range("L122").Value = 8550
$ws.Range("M122").Value = -6576.571599999999
$ws.Range("N122").Value = -13450
$ws.Range("H128").Value = 69000
$ws.Range("J128").Value = 69000
$ws.Range("L128").Value = 69000
$ws.Range("N128").Value = -78960

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 36.57895
$ws.Range("I38").Value = 10.714286
$ws.Range("J38").Value = 51.666668
$ws.Range("K38").Value = 32.142858
$ws.Range("L38").Value = 155.000004
$ws.Range("M38").Value = 314.857142
$ws.Range("N38").Value = -849.000004
$ws.Range("H98").Value = 388
$ws.Range("I98").Value = 180
$ws.Range("K98").Value = 540
$ws.Range("M98").Value = 958
$ws.Range("H117").Value = 557788
$ws.Range("J117").Value = 1666667
$ws.Range("L117").Value = 5000001
$ws.Range("N117").Value = -5006885
$ws.Range("H119").Value = 1399.75
$ws.Range("I119").Value = 1399.75
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 4199.25
$ws.Range("L119").Value = 0
$ws.Range("M119").Value = 638.75
$ws.Range("N119").ClearContents()
$ws.Range("H122").Value = 25092.334
$ws.Range("I122").Value = 37212
$ws.Range("J122").Value = 853
$ws.Range("K122").Value = 334908
$ws.Range("L122").Value = 7677
$ws.Range("M122").Value = -332458
$ws.Range("N122").Value = -12577
$ws.Range("H133").Value = 990
$ws.Range("I133").Value = 990
$ws.Range("K133").Value = 2970
$ws.Range("M133").Value = 2090
$ws.Range("H134").Value = 943
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H137").Value = 1181
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H139").Value = 5853.5713
$ws.Range("I139").Value = 2913.4285
$ws.Range("J139").Value = 8793.714
$ws.Range("K139").Value = 8740.2855
$ws.Range("L139").Value = 26381.142
$ws.Range("M139").Value = -3600.2855
$ws.Range("N139").Value = -36661.142
$ws.Range("H141").Value = 8999
$ws.Range("I141").Value = 8999
$ws.Range("K141").Value = 26997
$ws.Range("M141").Value = -21817

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 44444
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 44444
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 44444
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -44934
$ws.Range("H64").Value = 79082
$ws.Range("I64").Value = 77246
$ws.Range("K64").Value = 77246
$ws.Range("M64").Value = -76998
$ws.Range("H67").Value = 79082
$ws.Range("I67").Value = 77246
$ws.Range("K67").Value = 77246
$ws.Range("M67").Value = -76388
$ws.Range("H126").Value = 3127.3572
$ws.Range("I126").Value = 1976.9
$ws.Range("K126").Value = 5930.700000000001
$ws.Range("M126").Value = -3460.700000000001
$ws.Range("H132").Value = 6257.206
$ws.Range("I132").Value = 4093.353
$ws.Range("K132").Value = 12280.059
$ws.Range("M132").Value = -9750.059000000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4767416.5
$ws.Range("J46").Value = 5787.55
$ws.Range("L46").Value = 5787.55
$ws.Range("N46").Value = -6163.55
$ws.Range("H93").Value = 2700.1365
$ws.Range("I93").Value = 7233.8335
$ws.Range("K93").Value = 7233.8335
$ws.Range("M93").Value = -5985.8335
$ws.Range("H100").Value = 7356859
$ws.Range("I100").Value = 13161302
$ws.Range("K100").Value = 13161302
$ws.Range("M100").Value = -13160761
$ws.Range("H107").Value = 2400
$ws.Range("I107").Value = 2400
$ws.Range("K107").Value = 2400
$ws.Range("M107").Value = -480

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 7588.2354
$ws.Range("I15").Value = 7588.2354
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 7588.2354
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -7300.2354
$ws.Range("N15").ClearContents()
$ws.Range("H20").Value = 82998
$ws.Range("J20").Value = 82998
$ws.Range("L20").Value = 82998
$ws.Range("N20").Value = -83478
$ws.Range("H61").Value = 3254.0625
$ws.Range("I61").Value = 3357.6667
$ws.Range("K61").Value = 3357.6667
$ws.Range("M61").Value = -3065.6667
$ws.Range("H132").Value = 5290.4443
$ws.Range("I132").Value = 3725.1428
$ws.Range("J132").Value = 8173.8945
$ws.Range("K132").Value = 11175.4284
$ws.Range("L132").Value = 24521.6835
$ws.Range("M132").Value = -8645.428400000001
$ws.Range("N132").Value = -29581.6835
